# Update simulation output values (B2:E13) with unrounded floating point
# results, replacing the previously rounded integer values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2169.953739032937
$ws.Range("C2").Value = 2170.085420450942
$ws.Range("D2").Value = 2175.325749412518
$ws.Range("E2").Value = 2176.626862922192

$ws.Range("B3").Value = 2184.514313820853
$ws.Range("C3").Value = 2184.918646513586
$ws.Range("D3").Value = 2187.073464779258
$ws.Range("E3").Value = 2186.295744950701

$ws.Range("B4").Value = 2137.004174443122
$ws.Range("C4").Value = 2135.372170126485
$ws.Range("D4").Value = 2140.236464658321
$ws.Range("E4").Value = 2143.181371125894

$ws.Range("B5").Value = 2185.240087232406
$ws.Range("C5").Value = 2186.274067604017
$ws.Range("D5").Value = 2187.615832963864
$ws.Range("E5").Value = 2184.578900396438

$ws.Range("B6").Value = 2191.438427965602
$ws.Range("C6").Value = 2192.367416912643
$ws.Range("D6").Value = 2193.78184930315
$ws.Range("E6").Value = 2190.875211451854

$ws.Range("B7").Value = 2173.658829612546
$ws.Range("C7").Value = 2174.458556810208
$ws.Range("D7").Value = 2175.768832331416
$ws.Range("E7").Value = 2172.310782951094

$ws.Range("B8").Value = 2174.829099415711
$ws.Range("C8").Value = 2174.349416912381
$ws.Range("D8").Value = 2173.588269967609
$ws.Range("E8").Value = 2169.76364612606

$ws.Range("B9").Value = 2187.41046822575
$ws.Range("C9").Value = 2188.916812052127
$ws.Range("D9").Value = 2191.725020862515
$ws.Range("E9").Value = 2189.436098052507

$ws.Range("B10").Value = 1963.41524271665
$ws.Range("C10").Value = 1970.567967502973
$ws.Range("D10").Value = 2001.459579491783
$ws.Range("E10").Value = 2024.861244482418

$ws.Range("B11").Value = 1912.761164030131
$ws.Range("C11").Value = 1908.487280735242
$ws.Range("D11").Value = 1940.725141317559
$ws.Range("E11").Value = 1967.172106832768

$ws.Range("B12").Value = 1646.942019134947
$ws.Range("C12").Value = 1624.389466078262
$ws.Range("D12").Value = 1605.030983419213
$ws.Range("E12").Value = 1609.080127237416

$ws.Range("B13").Value = 1942.524931502256
$ws.Range("C13").Value = 1936.788313670181
$ws.Range("D13").Value = 1954.061397624516
$ws.Range("E13").Value = 1979.265623244521
